$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.306.62"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.601.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.51"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.97"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0856"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.59"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.612.58"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.294.85"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.07"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.90%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.97%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.18"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.48"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.45"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.446.46"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.97"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.568"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.80"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.738.66"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.919"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.759"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.83"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.38"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.48"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0500"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0950"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.21%  "
